# Decrement the "剩余" (remaining) value in column E by 1 for each data row,
# except row 36 whose value is left unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 - 1
    }
}
